$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A13 from "lab1" to "SA1"
$ws.Range("A13").Value = "SA1"

# Add new deadline value in B13, matching the text style used by B2:B12
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "September 09, 2024"

# Update the selected cell to match the new selection (D17)
$ws.Range("D17").Select()
